$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.517.13'
$ws.Range("E2").Value = '  +0.73%  '
$ws.Range("D3").Value = '1.728.51'
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.69'
$ws.Range("E5").Value = '  +2.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4789'
$ws.Range("E7").Value = '  +0.71%  '
$ws.Range("E8").Value = '  +1.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06224'
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("D10").Value = '1.731.34'
$ws.Range("E10").Value = '  +0.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07169'
$ws.Range("E11").Value = '  +1.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.66'
$ws.Range("E12").Value = '  +2.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6147'
$ws.Range("E13").Value = '  +4.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.529'
$ws.Range("E14").Value = '  +2.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.89'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").Value = '26.516.44'
$ws.Range("E17").Value = '  +0.69%  '
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006966'
$ws.Range("E19").Value = '  +2.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.63'
$ws.Range("E20").Value = '  +0.83%  '
$ws.Range("D21").Value = '1.952.46'
$ws.Range("E21").Value = '  +0.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.524'
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.914'
$ws.Range("E23").Value = '  +1.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.281'
$ws.Range("E24").Value = '  -0.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '136.45'
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("E26").Value = '  +0.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.790'
$ws.Range("E27").Value = '  +2.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.404'
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.61'
$ws.Range("E29").Value = '  -1.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.982'
$ws.Range("E30").Value = '  -0.25%  '
$ws.Range("E31").Value = '  +2.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.709'
$ws.Range("E32").Value = '  +0.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04589'
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.618'
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9941'
$ws.Range("E36").Value = '  +2.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6308'
$ws.Range("E37").Value = '  +1.97%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9156'
$ws.Range("E38").Value = '  -0.97%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.086'
$ws.Range("E39").Value = '  +9.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.406'
$ws.Range("E40").Value = '  -0.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '104.66'
$ws.Range("E41").Value = '  -6.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.003'
$ws.Range("E42").Value = '  +0.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01501'
$ws.Range("E43").Value = '  +2.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.581'
$ws.Range("E44").Value = '  +4.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3877'
$ws.Range("E45").Value = '  +1.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.978'
$ws.Range("E46").Value = '  +10.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1182'
$ws.Range("E47").Value = '  +1.63%  '
$ws.Range("E48").Value = '  +1.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '31.01'
$ws.Range("E49").Value = '  +2.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.802'
$ws.Range("E50").Value = '  +1.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.259'
$ws.Range("E51").Value = '  +3.91%  '
